$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated average_county_temperature (column K) values based on NOAA data,
# along with the dependent COP columns (R = worst_ashp_cop, S = best_ashp_cop)
# that were recalculated from the new temperature figures.

# Facility 1002037 (MO) - rows 5,6,7
$ws.Range("K5").Value = 17.71296296296294
$ws.Range("R5").Value = 1.969147761396439
$ws.Range("S5").Value = 2.170013198606573

$ws.Range("K6").Value = 17.71296296296294
$ws.Range("R6").Value = 1.872915723725898
$ws.Range("S6").Value = 2.048770944115581

$ws.Range("K7").Value = 17.71296296296294

# Facility 1002057 (GA) - rows 8,9
$ws.Range("K8").Value = -1.819444444444444
$ws.Range("R8").Value = 1.635040899548778
$ws.Range("S8").Value = 1.766799672330418

$ws.Range("K9").Value = -1.819444444444444
$ws.Range("R9").Value = 1.572614297115494
$ws.Range("S9").Value = 1.690895540926593

# Facility 1003352 (CO) - rows 26,27
$ws.Range("K26").Value = -0.763888888888889
$ws.Range("R26").Value = 1.650171687407155
$ws.Range("S26").Value = 1.78472092464565

$ws.Range("K27").Value = -0.763888888888889
$ws.Range("R27").Value = 1.586359976998275
$ws.Range("S27").Value = 1.707009404388715

# Facility 1003568 (CO) - rows 31,32
$ws.Range("K31").Value = -0.763888888888889
$ws.Range("R31").Value = 1.586359976998275
$ws.Range("S31").Value = 1.707009404388715

$ws.Range("K32").Value = -0.763888888888889
$ws.Range("R32").Value = 1.650171687407155
$ws.Range("S32").Value = 1.78472092464565

# Facility 1006904 (GA) - rows 41,42,43
$ws.Range("K41").Value = -1.819444444444444
$ws.Range("R41").Value = 1.635040899548778
$ws.Range("S41").Value = 1.766799672330418

$ws.Range("K42").Value = -1.819444444444444
$ws.Range("R42").Value = 1.572614297115494
$ws.Range("S42").Value = 1.690895540926593

$ws.Range("K43").Value = -1.819444444444444

$wb.Save()
